# Apply the edit described by the diff:
#  - Column D ("Produto") changes from "H2" to "H2V" for every data row (2-19)
#  - Column B ("Status") changes "Conceito" / "Contrução" -> "Estudo de Viabilidade"
#    and "DEMO" -> "Operando"
#  - Row 16, column E value changes from "0.1" to "0.3"
#  - Row 16, column C gets a (toggled) font style applied (Bold on/off), which
#    Excel records as a new, otherwise-identical cell style (applyFont="1")
#  - Active selection moves to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: "H2" -> "H2V" for all data rows ---
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "H2V"
}

# --- Column B: normalize status labels ---
$statusMap = @{
    4  = "Estudo de Viabilidade"
    5  = "Estudo de Viabilidade"
    6  = "Estudo de Viabilidade"
    7  = "Estudo de Viabilidade"
    10 = "Estudo de Viabilidade"
    12 = "Estudo de Viabilidade"
    13 = "Estudo de Viabilidade"
    14 = "Estudo de Viabilidade"
    15 = "Estudo de Viabilidade"
    16 = "Operando"
    17 = "Estudo de Viabilidade"
    18 = "Estudo de Viabilidade"
}

foreach ($r in $statusMap.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $statusMap[$r]
}

# --- Row 16, column E: "0.1" -> "0.3" ---
$ws.Cells.Item(16, 5).Value2 = "0.3"

# --- Row 16, column C: toggle bold to record the new (identical) style ---
$ws.Cells.Item(16, 3).Font.Bold = $true
$ws.Cells.Item(16, 3).Font.Bold = $false

# --- Move the active selection to D11 ---
$ws.Range("D11").Select()
